# "Added Week 15 simulations"
# Updates the Rushing and Receiving stat tables with the new weekly totals,
# renumbers the player-index column (A) so it stays contiguous, and adds a
# new Receiving row for D.Montgomery.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 (Z.Wilson)
$rushing.Range("C2").Value = 4
$rushing.Range("D2").Value = 2
$rushing.Range("E2").Value = 1
$rushing.Range("F2").Value = 5

# Row 7 (T.Johnson)
$rushing.Range("C7").Value = 26
$rushing.Range("D7").Value = 18

# Row 9 (M.White)
$rushing.Range("C9").Value = 7
$rushing.Range("D9").Value = 1

# Player-index column renumbering (close the gap left elsewhere in the roster)
$rushing.Range("A11").Value = 9
$rushing.Range("A12").Value = 10
$rushing.Range("A13").Value = 11

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Row 3 (T.Johnson)
$receiving.Range("C3").Value = 46
$receiving.Range("D3").Value = 29
$receiving.Range("G3").Value = 3
$receiving.Range("H3").Value = 2

# Row 6 (J.Crowder)
$receiving.Range("A6").Value = 4
$receiving.Range("C6").Value = 79
$receiving.Range("D6").Value = 60
$receiving.Range("E6").Value = 11

# Row 7 (K.Cole)
$receiving.Range("A7").Value = 5
$receiving.Range("C7").Value = 27
$receiving.Range("E7").Value = 12
$receiving.Range("F7").Value = 5
$receiving.Range("G7").Value = 6

# Row 8 (E.Moore) - index renumber only
$receiving.Range("A8").Value = 6

# Row 9 (B.Berrios)
$receiving.Range("A9").Value = 7
$receiving.Range("C9").Value = 39
$receiving.Range("D9").Value = 27
$receiving.Range("G9").Value = 4

# Row 10 (D.Mims)
$receiving.Range("A10").Value = 8
$receiving.Range("C10").Value = 10

# Row 11 (J.Smith) - index renumber only
$receiving.Range("A11").Value = 9

# Row 12 (N.Bawden) - index renumber only
$receiving.Range("A12").Value = 10

# Row 13 (T.Kroft)
$receiving.Range("A13").Value = 11
$receiving.Range("C13").Value = 6
$receiving.Range("D13").Value = 3
$receiving.Range("E13").Value = 0
$receiving.Range("F13").Value = 0
$receiving.Range("H13").Value = 0

# Row 14 (R.Griffin)
$receiving.Range("A14").Value = 12
$receiving.Range("C14").Value = 14
$receiving.Range("D14").Value = 9
$receiving.Range("E14").Value = 2
$receiving.Range("F14").Value = 1
$receiving.Range("G14").Value = 1
$receiving.Range("H14").Value = 1

# Row 15 (T.Wesco)
$receiving.Range("A15").Value = 13
$receiving.Range("C15").Value = 35
$receiving.Range("D15").Value = 22
$receiving.Range("E15").Value = 5
$receiving.Range("F15").Value = 3
$receiving.Range("G15").Value = 6
$receiving.Range("H15").Value = 3

# New row 16 - D.Montgomery. Copy the formatting from the row above first
# (bold/bordered index style on column A) and then fill in the values.
$receiving.Range("A15").Copy()
$receiving.Range("A16").PasteSpecial(-4122)

$receiving.Range("A16").Value = 14
$receiving.Range("B16").Value = "D.Montgomery"
$receiving.Range("C16").Value = 4
$receiving.Range("D16").Value = 2
$receiving.Range("E16").Value = 0
$receiving.Range("F16").Value = 0
$receiving.Range("G16").Value = 0
$receiving.Range("H16").Value = 0
